$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the "Diaria" table (rows 87-95), matching the
# additional daily series published for late Sept / early Oct 2021.
$rows = @(
    @{ Row = 87; A = "21-09-2021"; B = 50000 }
    @{ Row = 88; A = "22-09-2021"; B = 50000 }
    @{ Row = 89; A = "23-09-2021"; B = 50000 }
    @{ Row = 90; A = "24-09-2021"; B = 50000; C = 110000; E = 25000; F = 0;    G = 2.18 }
    @{ Row = 91; A = "27-09-2021"; B = 50000; C = 110000; E = 20000; F = 5000; G = 2.21 }
    @{ Row = 92; A = "28-09-2021"; B = 50000; C = 80000;  E = 25000; F = 0;    G = 2.2 }
    @{ Row = 93; A = "29-09-2021"; B = 50000 }
    @{ Row = 94; A = "30-09-2021"; B = 50000 }
    @{ Row = 95; A = "01-10-2021"; B = 50000 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B

    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($rowNum, 3).Value = $r.C
    }

    # Column D always carries a value: adjudicado total minus demandado when
    # there was an auction that week, otherwise 0.
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($rowNum, 4).Value = 25000
    } else {
        $ws.Cells.Item($rowNum, 4).Value = 0
    }

    if ($r.ContainsKey("E")) {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    }
    if ($r.ContainsKey("F")) {
        $ws.Cells.Item($rowNum, 6).Value = $r.F
    }
    if ($r.ContainsKey("G")) {
        $ws.Cells.Item($rowNum, 7).Value = $r.G
    }
}
